$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5 ("Crumpet GEF" facility) needs the same cell formatting (styles) as the
# existing data rows above it (row 4), but only for columns A:J - row 4's K column
# (payment exchange rate) formatting/value must NOT be copied down, since the new
# rows don't carry a value in that column.
$ws.Range("A4:J4").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new fixture rows used by the rounding-error e2e test.
# Values are entered column-by-column (matching how the shared-string table
# ends up ordered in the saved workbook).
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("A6").Value = "Scone GEF"

$ws.Range("B5").Value = 20001371
$ws.Range("B6").Value = 20001371

$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("C6").Value = "Scone exporter"

$ws.Range("D5").Value = "GBP"
$ws.Range("D6").Value = "GBP"

$ws.Range("E5").Value = 7000000
$ws.Range("E6").Value = 770000

$ws.Range("F5").Value = 3938753.8
$ws.Range("F6").Value = 761579.37

$ws.Range("G5").Value = 777
$ws.Range("G6").Value = 777

$ws.Range("H5").Value = 456
$ws.Range("H6").Value = 456.77

$ws.Range("I5").Value = "GBP"
$ws.Range("I6").Value = "GBP"

$ws.Range("J5").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Reflect the author's final cell selection / scroll position before saving.
$ws.Range("F8").Select() | Out-Null
